$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (trial counts) changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) values updated
$ws.Range("B2").Value = 50.720884012982765
$ws.Range("C2").Value = 56.772752221374418
$ws.Range("D2").Value = 53.094941844320054
$ws.Range("E2").Value = 60.124740101386323

# Row 3 (STR) values updated
$ws.Range("B3").Value = 45.78411335805194
$ws.Range("C3").Value = 51.026369555577283
$ws.Range("D3").Value = 51.476800363831472
$ws.Range("E3").Value = 56.726197686557249

# Selection narrowed from B1:AY3 to B1:E3
$ws.Range("B1:E3").Select()
